# "save tex - Kannai updates"
#
# 1) The "datetimeFigureOut" date placeholder that lives on the slide
#    master and on every slide layout was refreshed from 2022. 02. 28.
#    to 2022. 04. 24.
# 2) A typo was fixed on slide 1: "meghívot" -> "meghívott".

$p = $ppt.ActivePresentation
$newDate = "2022. 04. 24."

# --- 1. Update the date placeholder text everywhere it appears -------

# Slide master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Fix the typo on slide 1 ---------------------------------------

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "Az eredetileg meghívot függvény ezt az értéket adja vissza") {
            $shp.TextFrame.TextRange.Text = "Az eredetileg meghívott függvény ezt az értéket adja vissza"
        }
    }
}
